# Generate Report for Handback
# Refresh the two sample UUID-named files and their associated
# handoff/handback timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newUuid1 = "2639c563-4a9d-40ff-a374-978388efcecf"
$newUuid2 = "ffffc6449859-46b8-4ccd-8643-633cf92b4108"

$newXlfZh = "$newUuid1.1240c37d6784f132086865b513d001413ae18443.zh-cn.xlf"
$newXlfDe = "$newUuid1.1240c37d6784f132086865b513d001413ae18443.de-de.xlf"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("G2").Value = "2016-08-12 03:23:38"
$wsOverview.Range("G3").Value = "2016-08-12 03:23:38"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newUuid1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newUuid2.md"
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newUuid1.md"
$wsZh.Range("I2").Value = "$newUuid1.md"
$wsZh.Range("A3").Value = "$newUuid2.md"
$wsZh.Range("I3").Value = "$newUuid2.md"

$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("J2").Value = $newXlfZh
$wsZh.Range("H2").Value = "2016-08-12 03:23:33"
$wsZh.Range("K2").Value = "2016-08-12 03:23:48"

$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("J3").Value = $newXlfZh
$wsZh.Range("H3").Value = "2016-08-12 03:23:33"
$wsZh.Range("K3").Value = "2016-08-12 03:23:48"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    }
}

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newUuid1.md"
$wsDe.Range("I2").Value = "$newUuid1.md"
$wsDe.Range("A3").Value = "$newUuid2.md"
$wsDe.Range("I3").Value = "$newUuid2.md"

$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("H2").Value = "2016-08-12 03:23:38"
$wsDe.Range("K2").Value = "2016-08-12 03:23:56"

$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("H3").Value = "2016-08-12 03:23:38"
$wsDe.Range("K3").Value = "2016-08-12 03:23:56"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newUuid2.md"
    }
}
